# Changed date of sample flight
# The sample flight log on Sheet1 had its Start/Completion time values
# shifted from 2019-01-01 to 2019-01-10 (same time-of-day, 9 days later).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Flight 1): Start time / Completion time
$ws.Range("B2").Value = 43475.647037037037
$ws.Range("C2").Value = 43475.647893518515

# Row 3 (Flight 2): Start time / Completion time
$ws.Range("B3").Value = 43475.647974537038
$ws.Range("C3").Value = 43475.648449074077

# Row 4 (Flight 3): Start time / Completion time
$ws.Range("B4").Value = 43475.648472222223
$ws.Range("C4").Value = 43475.648877314816

# The active selection on the sheet moved from C13 to C14
$ws.Range("C14").Select()
